# Append a new row of scraped data to the "ランサーズ" sheet and refresh the
# "取得日時" (fetched-at) timestamp on every existing row to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-01-26 18:32:11"

# Refresh the fetch timestamp column (A) for all existing data rows (2-9).
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# This listing's price/terms text changed between the two scrapes.
$ws.Range("D8").Value = "5,000 円 ~ 6,000 円 / 募集期間 5 日、取引期間 0 日"

# Column D (価格) needed to widen to fit the longer text above.
$ws.Columns.Item(4).ColumnWidth = 38.166666666666664

# Append the newly scraped listing as row 10.
$ws.Range("A10").Value = $newTimestamp
$ws.Range("B10").Value = "【急募】NASのショートカットリンク一括編集PowerShellコード作成"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5479430"
$ws.Range("G10").Value = 10

$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5479430")
$ws.Range("F10").Style = "Hyperlink"
